$p = $ppt.ActivePresentation

# 1. Delete slide 10 (id=262, the "Process" slide with Lexer/Generator bullet
#    text) which is being dropped from the deck. The remaining "Process"
#    slide (id=266, the detailed flow-chart version) shifts up to take its
#    place (position 10).
$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 262) {
        $target = $p.Slides.Item($i)
        break
    }
}
if ($target -ne $null) {
    $target.Delete()
}

# 2. Update the cached "datetimeFigureOut" date field text from 2019-06-19
#    to 2019-06-20 on the slide master and every slide layout.
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

foreach ($shp in $master.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "2019-06-19") {
            $shp.TextFrame.TextRange.Text = "2019-06-20"
        }
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    foreach ($shp in $layout.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "2019-06-19") {
                $shp.TextFrame.TextRange.Text = "2019-06-20"
            }
        }
    }
}
